$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap cell values: A2 becomes the buurtcode "BU02220303", B2 becomes "Unknown"
$ws.Range("A2").Value = "BU02220303"
$ws.Range("B2").Value = "Unknown"

# Remove the bold/custom font formatting previously applied to A2
$ws.Range("A2").ClearFormats()

# Column A width (~17.57 chars wide, matching the manually-resized column
# from the authored workbook; COM's character-unit ColumnWidth quantizes to
# the nearest 1/6 character so 16.71 is the closest achievable setting)
$ws.Columns.Item(1).ColumnWidth = 16.71

# Selection moves to A3
$ws.Range("A3").Select() | Out-Null
